{"js": "const tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// 1) Update the three data-cell values in the first (remaining) table.\n//    getCell uses 0-based (row, column) indices, so row 1 / cols 1-3\n//    are the \"Nom de l'etape\" / \"Description de l'\u00e9tape\" / \"R\u00e9sultat Attendu\"\n//    data cells of the second table row.\nconst t1 = tables.items[0];\nt1.getCell(1, 1).value = \"step1\";\nt1.getCell(1, 2).value = \"lkfdf\";\nt1.getCell(1, 3).value = \"hyjyhtg\";\n\n// 2) Remove the whole second table (the duplicate header-only table).\nconst t2 = tables.items[1];\nt2.delete();\n\n// 3) Remove the \"fg\" Heading2 paragraph that introduced the now-deleted table.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet fgParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"fg\") {\n    fgParagraph = paragraphs.items[i];\n    break;\n  }\n}\nif (fgParagraph) {\n  fgParagraph.delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the three data-cell values in the first (remaining) table.\n$t1 = $d.Tables.Item(1)\n$t1.Cell(2, 2).Range.Text = \"step1\"\n$t1.Cell(2, 3).Range.Text = \"lkfdf\"\n$t1.Cell(2, 4).Range.Text = \"hyjyhtg\"\n\n# 2) Remove the whole second table (the duplicate \"Index/Nom/Description/Resultat\" header table).\n$t2 = $d.Tables.Item(2)\n$t2.Delete()\n\n# 3) Remove the \"fg\" Heading2 paragraph that introduced the now-deleted table.\n$rng = $d.Content\n$rng.Find.Text = \"fg\"\n$rng.Find.Execute() | Out-Null\n$p = $rng.Paragraphs.Item(1)\n$fullPara = $d.Range($p.Range.Start, $p.Range.End + 1)\n$fullPara.Delete()\n"}
